# Release v0.1.0-beta: Fix validation errors and update canonical URL
#
# Updates the "Metadata" sheet (Version, Status, Date, Description) and the
# "Elements" sheet (Definition of the root Extension row) to match the
# published v0.1.0-beta content.

$wb = $excel.ActiveWorkbook

$description = "Extension to link nursing interventions to the patient goals they are intended to achieve. Supports goal-directed care planning and intervention tracking."

# --- Metadata sheet ---
$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B3").Value = "0.1.0"              # Version: 1.0.0 -> 0.1.0
$meta.Range("B6").Value = "draft"               # Status: active -> draft
$meta.Range("B8").Value = "2025-12-26T14:13:58+00:00"  # Date
$meta.Range("B11").Value = $description         # Description: (blank) -> text

# --- Elements sheet ---
$elements = $wb.Worksheets.Item("Elements")

# Definition column (M) for the root "Extension" element row (row 2)
$elements.Range("M2").Value = $description
